$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Key to Variables")
$ws.Activate()

# ---------------------------------------------------------------------------
# The input data for "Share of Electricity Imports and Exports by Hour"
# (acronym SoEIaEbH) is being replaced by a new, endogenously-calculated
# variable "Max and Min Hourly Electriicty Imports and Exports" (acronym
# MaMHEIaE). In the "Key to Variables" sheet, this means:
#   1) Insert a new row (it lands right after the "LFHVM" row, i.e. row 100)
#      describing the new MaMHEIaE variable.
#   2) Remove the old row describing SoEIaEbH further down the sheet.
# ---------------------------------------------------------------------------

# 1) Insert the new "MaMHEIaE" row just below row 99 ("LFHVM"), so it becomes
#    row 100, copying the Top Level Folder + Importance values from the row
#    above it (both are "elec" / "high").
$ws.Rows.Item(100).Insert()
$ws.Range("A100").Value2 = $ws.Range("A99").Value2
$ws.Range("B100").Value2 = "MaMHEIaE"
$ws.Range("C100").Value2 = "Max and Min Hourly Electriicty Imports and Exports"
$ws.Range("F100").Value2 = $ws.Range("F99").Value2

# 2) Delete the old "SoEIaEbH" row. After the insert above, it has shifted
#    down from row 112 to row 113.
$ws.Rows.Item(113).Delete()

# ---------------------------------------------------------------------------
# Update the sheet's view state (frozen pane / selection) to reflect where
# the author was working when they made the edit.
# ---------------------------------------------------------------------------
$ws.Range("A2").Select()
$excel.ActiveWindow.FreezePanes = $true
$ws.Range("H100").Select()
